# LeilaVerbeek_0940566_Centralebank_V1 - "reupload centrale bank v1 / een onderdeel
# te veel bij A5" - the quality-requirement sub-slide about "maakbaarheid"
# (makeability) is a duplicate/unwanted item under the "Kwaliteitseisen
# eindresultaat" (A5) section, so it is removed entirely, and the matching
# bullet is dropped from the table-of-contents slide.

$p = $ppt.ActivePresentation

# --- 1. Table of contents slide ("Inhoud"): drop the "Maakbaarheid" bullet ---
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    foreach ($shape in $slide.Shapes) {
        if (-not $shape.HasTextFrame) { continue }
        if (-not $shape.TextFrame.HasText) { continue }
        $tr = $shape.TextFrame.TextRange
        $count = $tr.Paragraphs().Count
        for ($i = 1; $i -le $count; $i++) {
            $paragraph = $tr.Paragraphs($i, 1)
            if ($paragraph.Text.Trim() -eq "Maakbaarheid") {
                $paragraph.Delete()
                # Re-fit the text now that there is one fewer line (was
                # normAutofit with a fontScale/lnSpcReduction shrink; after
                # removing a bullet it fits natively).
                $shape.TextFrame.AutoSize = 2
                break
            }
        }
    }
}

# --- 2. Delete the "maakbaarheid" quality-requirement slide ---
# (title "Analyseren" / "Kwaliteitseisen eindresultaat ... maakbaarheid")
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $slideText = ""
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.HasText) {
                $slideText = $slideText + $shape.TextFrame.TextRange.Text
            }
        }
    }
    if ($slideText -like "*maakbaarheid*") {
        $slide.Delete()
        break
    }
}
